$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the header style (bold/bordered/centered) from G1 onto the new H1 header cell
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("H1").Value = "Save"
$ws.Range("H2").Value = 0
